$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1: next week's date, written as text (not coerced to a date
# serial) and formatted the same as the other header cells (C1).
$ws.Range("D1").NumberFormat = "@"
$ws.Range("D1").Value = "2025-02-06"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# New attendance-prediction cell D2, matching the plain (unstyled) look of C2.
$ws.Range("D2").Value = "P"
